# Applies the cibmtr-reporting-ig update to the ValueSet metadata workbook.
#
# Sheet "Metadata" (sheet1): a new "Jurisdiction" row is inserted after the
# two "Contact" rows, shifting Description/Purpose/Copyright/Immutable down
# by one row (Immutable becomes a brand-new row 16). Several values are
# also updated in place (Version, Status, Date, the two Contact values).
#
# Sheet "Include from SNOMED CT" (sheet2) is unaffected in terms of actual
# cell content; the diff there is purely a consequence of shared-string
# renumbering caused by the Metadata sheet edits.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# -- Simple in-place value updates -----------------------------------------
$ws.Range("B3").Value = "0.1.7"
$ws.Range("B6").Value = "draft"
$ws.Range("B8").Value = "2024-08-27T12:23:18-05:00"
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# -- Make room for the new "Jurisdiction" row by pushing rows 12-15 down ---
# First, grow the format of row 16 from row 15 (same look as the rest of
# the property table) before we start moving content into it.
$ws.Range("A15:B15").Copy()
$ws.Range("A16:B16").PasteSpecial(-4122)

# Shift the existing Description/Purpose/Copyright/Immutable rows down by
# one (row 15 <- row14, row14 <- row13, row13 <- row12), from the bottom up
# so we never overwrite a value before it has been copied forward.
$ws.Range("A15").Value = $ws.Range("A14").Value()
$ws.Range("B15").Value = $ws.Range("B14").Value()

$ws.Range("A14").Value = $ws.Range("A13").Value()
$ws.Range("B14").Value = $ws.Range("B13").Value()

$ws.Range("A13").Value = $ws.Range("A12").Value()
$ws.Range("B13").Value = $ws.Range("B12").Value()

# Row 12 becomes the new Jurisdiction row (value left blank).
$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""

# Row 16 is the new home for "Immutable" / "BooleanType[null]".
$ws.Range("A16").Value = "Immutable"
$ws.Range("B16").Value = "BooleanType[null]"
